$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 760 (pushes the old row 760..801 down to
# 761..802, matching the diff's "row inserted, everything below shifts
# by one" shape). Excel semantics: Rows.Item(760).Insert() shifts the
# existing row 760 (and below) down by one and leaves a blank row 760.
$ws.Rows.Item(760).Insert()

# The date column must stay plain text ("2026/02/01"), not be
# reinterpreted as a date serial number. Assigning it as a formula that
# evaluates to a text string ("2026/02/01") sidesteps Excel's "looks
# like a date" auto-conversion that a plain .Value assignment would
# trigger; copy/paste-values then bakes it into a literal text value
# without leaving any NumberFormat/style override on the cell.
$ws.Cells.Item(760, 1).Formula = '="2026/02/01"'
$ws.Cells.Item(760, 1).Copy()
$ws.Cells.Item(760, 1).PasteSpecial(-4163)  # xlPasteValues

$ws.Cells.Item(760, 2).Value = "日"
$ws.Cells.Item(760, 3).Value = 7
$ws.Cells.Item(760, 4).Value = 201
